$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '74.857.09'
$ws.Range("E2").Value = '  +1.85%  '

# Row 3
$ws.Range("D3").Value = '2.813.25'
$ws.Range("E3").Value = '  +7.64%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.89'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.38'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.191'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.07%  '

# Row 10
$ws.Range("D10").Value = '2.810.38'
$ws.Range("E10").Value = '  +7.59%  '

# Row 11
$ws.Range("E11").Value = '  -1.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.370'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.31%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.83'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.20%  '

# Row 14
$ws.Range("D14").Value = '3.330.95'
$ws.Range("E14").Value = '  +7.62%  '

# Row 15
$ws.Range("D15").Value = '74.784.08'
$ws.Range("E15").Value = '  +1.62%  '

# Row 16
$ws.Range("E16").Value = '  -0.17%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.79'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.05%  '

# Row 18
$ws.Range("D18").Value = '2.807.38'
$ws.Range("E18").Value = '  +7.00%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.92'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.06%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.27'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.86%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.23'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.39%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.44%  '

# Row 24
$ws.Range("E24").Value = '  -0.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.53%  '

# Row 26
$ws.Range("D26").Value = '2.946.78'
$ws.Range("E26").Value = '  +7.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.57%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.63'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.40%  '

# Row 29
$ws.Range("E29").Value = '  +11.21%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.14%  '

# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.29%  '

# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '513.38'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.69'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.49%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.02'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.90'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.118'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.39'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.62%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '185.35'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +14.68%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.09%  '

# Row 43
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.338'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.41%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.66'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.20'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.44%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.96'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.33'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0849'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.58%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.570'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +8.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.69'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.35%  '

# Row 51
$ws.Range("E51").Value = '  +8.64%  '
